$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 80
$ws.Range("I2").Value = 190
$ws.Range("J2").Value = 810
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 216
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = 137
$ws.Range("O2").Value = 1
$ws.Range("R2").Value = 14
$ws.Range("S2").Value = 76
$ws.Range("T2").Value = 125
$ws.Range("U2").Value = 13
$ws.Range("V2").Value = 1274
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1278
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 9
